$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.703.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.523.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0805"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.919.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.473.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.557.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "286.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.17%  "

$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0773"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "120.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0302"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.003.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
